# Updated cryptos list on Wed Oct 16 19:01:39 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. '601.40', '67.825.85').
# Prefix with a literal apostrophe so Excel stores them as text instead of
# silently coercing them into numbers and dropping formatting like trailing
# zeros or thousands-separator dots.

$ws.Range('D2').Value = "'" + '67.825.85'
$ws.Range('D3').Value = "'" + '2.619.02'
$ws.Range('E3').Value = '  +0.95%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'" + '601.40'
$ws.Range('E5').Value = '  +1.21%  '
$ws.Range('D6').Value = "'" + '154.51'
$ws.Range('E6').Value = '  -0.40%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +0.86%  '
$ws.Range('D9').Value = "'" + '2.616.80'
$ws.Range('E9').Value = '  +0.91%  '
$ws.Range('D10').Value = "'" + '0.126'
$ws.Range('E10').Value = '  +9.30%  '
$ws.Range('E11').Value = '  +0.93%  '
$ws.Range('D12').Value = "'" + '5.25'
$ws.Range('E12').Value = '  +0.58%  '
$ws.Range('D13').Value = "'" + '0.354'
$ws.Range('E13').Value = '  -1.63%  '
$ws.Range('D14').Value = "'" + '27.64'
$ws.Range('D15').Value = "'" + '0.0000186'
$ws.Range('E15').Value = '  +2.57%  '
$ws.Range('D16').Value = "'" + '3.099.86'
$ws.Range('E16').Value = '  +1.46%  '
$ws.Range('D17').Value = "'" + '67.683.85'
$ws.Range('E17').Value = '  +1.01%  '
$ws.Range('D18').Value = "'" + '2.617.04'
$ws.Range('E18').Value = '  +0.89%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = "'" + '366.65'
$ws.Range('E19').Value = '  +3.00%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').Value = "'" + '11.18'
$ws.Range('E20').Value = '  -1.24%  '
$ws.Range('D21').Value = "'" + '7.64'
$ws.Range('E21').Value = '  -2.44%  '
$ws.Range('E22').Value = '  -0.70%  '
$ws.Range('E23').Value = '  -2.59%  '
$ws.Range('D24').Value = "'" + '0.999'
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('D25').Value = "'" + '9.85'
$ws.Range('E25').Value = '  -7.15%  '
$ws.Range('D26').Value = "'" + '67.03'
$ws.Range('E26').Value = '  -0.28%  '
$ws.Range('D27').Value = "'" + '2.747.28'
$ws.Range('D28').Value = "'" + '0.0000103'
$ws.Range('E28').Value = '  -0.68%  '
$ws.Range('D29').Value = "'" + '575.25'
$ws.Range('E29').Value = '  -4.90%  '
$ws.Range('D30').Value = "'" + '1.00'
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('E31').Value = '  -3.00%  '
$ws.Range('E32').Value = '  -2.46%  '
$ws.Range('E33').Value = '  +0.46%  '
$ws.Range('D34').Value = "'" + '0.132'
$ws.Range('E34').Value = '  -1.95%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('E36').Value = '  -3.99%  '
$ws.Range('D37').Value = "'" + '4.92'
$ws.Range('E37').Value = '  -2.40%  '
$ws.Range('D38').Value = "'" + '159.01'
$ws.Range('E38').Value = '  +3.11%  '
$ws.Range('D39').Value = "'" + '19.35'
$ws.Range('E39').Value = '  +0.23%  '
$ws.Range('E40').Value = '  -0.36%  '
$ws.Range('D41').Value = "'" + '5.34'
$ws.Range('E41').Value = '  -2.88%  '
$ws.Range('E42').Value = '  +1.77%  '
$ws.Range('D43').Value = "'" + '2.57'
$ws.Range('E43').Value = '  -3.88%  '
$ws.Range('D44').Value = "'" + '41.17'
$ws.Range('E44').Value = '  -1.02%  '
$ws.Range('D45').Value = "'" + '0.999'
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('E46').Value = '  -0.12%  '
$ws.Range('D47').Value = "'" + '156.98'
$ws.Range('E47').Value = '  +0.26%  '
$ws.Range('D48').Value = "'" + '0.0₆0287'
$ws.Range('E48').Value = '  -7.82%  '
$ws.Range('E49').Value = '  -0.57%  '
$ws.Range('D50').Value = "'" + '21.00'
$ws.Range('E50').Value = '  -2.42%  '
$ws.Range('D51').Value = "'" + '0.623'
$ws.Range('E51').Value = '  +1.11%  '
